$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append a new list item "At least 10 questions" after the last
#    paragraph ("(I am doing 4.5/Binary)"), inheriting the same
#    ListParagraph style / numbering / run formatting.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.InsertBefore("At least 10 questions")

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the last paragraph
#    to the end of the "Would want to be able to win" paragraph.
#
#    The engine mis-places a bookmark that is added at a position
#    exactly on a paragraph mark, so work around it by temporarily
#    inserting a placeholder character at the target spot, anchoring
#    the bookmark next to that placeholder (a normal text position),
#    and then deleting the placeholder again. Re-adding a bookmark
#    with the existing name "_GoBack" also removes the old occurrence
#    automatically.
# ------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Would want to be able to win") {
        $targetPara = $cand
        break
    }
}

$endPos = $targetPara.Range.End - 1
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Delete()
